$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("section" and everything after it
# shifts one column to the right, from D.. to E..).
$ws.Columns("D:D").Insert()

# New "course" column header + its instructional sample value.
$ws.Range("D1").Value = "course"
$ws.Range("D2").Value = "1,2,3 (depending on MBA,BBA,Btech respectively)"

# The email hyperlink cell shifted from I2 to J2 with the column insert;
# recreate the mailto hyperlink there (and restore its text style).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:777@dd.com")
$ws.Range("J2").NumberFormat = "@"

# New column gets its own best-fit width (matches the width Excel computed
# for the "course" header / instructional text).
$ws.Range("D1").ColumnWidth = 45.17

# Selection ends up on D6 after the edit.
$ws.Range("D6").Select()
